$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7939458
$ws.Range("I76").Value = 10103737
$ws.Range("J76").Value = 3766.6667
$ws.Range("K76").Value = 10103737
$ws.Range("L76").Value = 3766.6667
$ws.Range("M76").Value = -10103422
$ws.Range("N76").Value = -4396.6667
$ws.Range("H79").Value = 7939458
$ws.Range("I79").Value = 10103737
$ws.Range("J79").Value = 3766.6667
$ws.Range("K79").Value = 10103737
$ws.Range("L79").Value = 3766.6667
$ws.Range("M79").Value = -10102645
$ws.Range("N79").Value = -5950.6667
$ws.Range("H98").Value = 467707.75
$ws.Range("I98").Value = 534123.1
$ws.Range("J98").Value = 2800
$ws.Range("K98").Value = 534123.1
$ws.Range("L98").Value = 2800
$ws.Range("M98").Value = -532625.1
$ws.Range("N98").Value = -5796
$ws.Range("H122").Value = 467707.75
$ws.Range("I122").Value = 534123.1
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 1602369.3
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -1599919.3
$ws.Range("N122").Value = -13300
$ws.Range("H133").Value = 49999.855
$ws.Range("J133").Value = 49999.855
$ws.Range("L133").Value = 49999.855
$ws.Range("N133").Value = -60119.855
$ws.Range("H137").Value = 125002620
$ws.Range("I137").Value = 166668340
$ws.Range("J137").Value = 5501.5
$ws.Range("K137").Value = 500005020
$ws.Range("L137").Value = 16504.5
$ws.Range("M137").Value = -500002470
$ws.Range("N137").Value = -21604.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 23900
$ws.Range("I63").Value = 38000
$ws.Range("K63").Value = 38000
$ws.Range("M63").Value = -37314
$ws.Range("H66").Value = 23900
$ws.Range("I66").Value = 38000
$ws.Range("K66").Value = 190000
$ws.Range("M66").Value = -186568
$ws.Range("H102").Value = 2986.85
$ws.Range("I102").Value = 3045.5
$ws.Range("K102").Value = 3045.5
$ws.Range("M102").Value = -1423.5
$ws.Range("H122").Value = 2729.5
$ws.Range("I122").Value = 2699.4443
$ws.Range("K122").Value = 8098.3329
$ws.Range("M122").Value = -5648.3329

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 624.57574
$ws.Range("I107").Value = 628.46875
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 628.46875
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1291.53125
$ws.Range("N107").Value = -4340

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6246.982
$ws.Range("I31").Value = 5716
$ws.Range("J31").Value = 6324.4165
$ws.Range("K31").Value = 5716
$ws.Range("L31").Value = 6324.4165
$ws.Range("N31").Value = -6914.4165
$ws.Range("M31").Value = -5421
$ws.Range("H34").Value = 6246.982
$ws.Range("I34").Value = 5716
$ws.Range("J34").Value = 6324.4165
$ws.Range("K34").Value = 5716
$ws.Range("L34").Value = 6324.4165
$ws.Range("N34").Value = -6728.4165
$ws.Range("M34").Value = -5514
$ws.Range("H122").Value = 1423.8
$ws.Range("I122").Value = 1453.0714
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 4359.2142
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = -1909.2142
$ws.Range("N122").Value = -7942

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3135486.2
$ws.Range("I113").Value = 333
$ws.Range("J113").Value = 5051413
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 15154239
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -15158579
$ws.Range("H117").Value = 534.7143
$ws.Range("I117").Value = 208
$ws.Range("J117").Value = 970.3333
$ws.Range("K117").Value = 624
$ws.Range("L117").Value = 2910.9999
$ws.Range("M117").Value = 2818
$ws.Range("N117").Value = -9794.999899999999
$ws.Range("H131").Value = 6668218
$ws.Range("J131").Value = 7409070.5
$ws.Range("L131").Value = 22227211.5
$ws.Range("N131").Value = -22237291.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H80").Value = 55558280
$ws.Range("J80").Value = 166669440
$ws.Range("L80").Value = 166669440
$ws.Range("N80").Value = -166671436
$ws.Range("H83").Value = 55558280
$ws.Range("J83").Value = 166669440
$ws.Range("L83").Value = 833347200
$ws.Range("N83").Value = -833357184
$ws.Range("H102").Value = 3610.647
$ws.Range("I102").Value = 3607.3635
$ws.Range("J102").Value = 3616.6667
$ws.Range("K102").Value = 3607.3635
$ws.Range("L102").Value = 3616.6667
$ws.Range("M102").Value = -1985.3635
$ws.Range("N102").Value = -6860.6667
$ws.Range("H122").Value = 2792.423
$ws.Range("I122").Value = 2905.348
$ws.Range("J122").Value = 1926.6666
$ws.Range("K122").Value = 8716.044
$ws.Range("L122").Value = 5779.9998
$ws.Range("M122").Value = -6266.044
$ws.Range("N122").Value = -10679.9998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3997.8572
$ws.Range("J122").Value = 3997.8572
$ws.Range("L122").Value = 11993.5716
$ws.Range("N122").Value = -16893.5716
$ws.Range("H124").Value = 36500
$ws.Range("J124").Value = 36500
$ws.Range("L124").Value = 36500
$ws.Range("M124").Value = -46320
$ws.Range("H136").Value = 8617.467000000001
$ws.Range("I136").Value = 4076.1667
$ws.Range("J136").Value = 11645
$ws.Range("K136").Value = 12228.5001
$ws.Range("L136").Value = 34935
$ws.Range("M136").Value = -9678.500100000001
$ws.Range("N136").Value = -40035

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 77534.46000000001
$ws.Range("I126").Value = 91386.17999999999
$ws.Range("J126").Value = 1350
$ws.Range("K126").Value = 274158.54
$ws.Range("L126").Value = 4050
$ws.Range("M126").Value = -271688.54
$ws.Range("N126").Value = -8990
$ws.Range("H132").Value = 4069.8572
$ws.Range("I132").Value = 4161.0415
$ws.Range("J132").Value = 3870.9092
$ws.Range("K132").Value = 12483.1245
$ws.Range("L132").Value = 11612.7276
$ws.Range("M132").Value = -9953.124500000002
$ws.Range("N132").Value = -16672.7276
$ws.Range("H136").Value = 2733.4707
$ws.Range("I136").Value = 867.9
$ws.Range("J136").Value = 5398.5713
$ws.Range("K136").Value = 2603.7
$ws.Range("L136").Value = 16195.7139
$ws.Range("M136").Value = -53.69999999999982
$ws.Range("N136").Value = -21295.7139
